$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old column J ("compSel" header / value 1) is removed entirely.
$ws.Range("J1:J2").EntireColumn.Delete()

# Three new columns are inserted before the (old) column E, pushing the
# remaining contact-us columns to the right. Excel inherits the left
# neighbour's (column D) formatting for the freshly inserted cells, which
# already matches the target styling for row 1 (s=2) and for D2/E2/G2 (s=8).
$ws.Range("E1:G1").EntireColumn.Insert()

# New header row values (row 1) - keep inherited style (s=2)
$ws.Range("E1").Value = "NegZipCode1"
$ws.Range("F1").Value = "NegZipCode2"
$ws.Range("G1").Value = "PosZipCode"

# New data row values (row 2).
# A leading apostrophe forces text + keeps the cell's quote-prefix flag set,
# matching the inherited style (s=8) for E2 and G2.
$ws.Range("E2").Value = "'abcd"

# F2 ends up with a brand-new style (default font, quote-prefixed) rather
# than inheriting column D's Consolas styling, so clear its formatting
# before writing the quote-prefixed value.
$ws.Range("F2").ClearFormats()
$ws.Range("F2").Value = "'1$#abcd234"

$ws.Range("G2").Value = "'10002"

# Match the saved selection (entire column M, active cell M1).
$ws.Columns("M").Select() | Out-Null
